$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder country labels (sorted rank changed after the data refresh) ---
# Camerun moved up three positions (past Dinamarca, Corea del Sur, Marruecos).
# Use a temp placeholder to avoid clobbering "Camerun" before it is copied out.
$ws.Cells.Item(65,1).Value = "TEMP_PLACEHOLDER_1"
$ws.Cells.Item(66,1).Value = "Marruecos"
$ws.Cells.Item(67,1).Value = "Corea del Sur"
$ws.Cells.Item(68,1).Value = "Dinamarca"
$ws.Cells.Item(65,1).Value = "Camerun"

# Islas Malvinas moved up one position, swapping with Groenlandia.
$ws.Cells.Item(209,1).Value = "TEMP_PLACEHOLDER_2"
$ws.Cells.Item(210,1).Value = "Groenlandia"
$ws.Cells.Item(209,1).Value = "Islas Malvinas"

# --- Update the "data refreshed at" timestamp ---
$ws.Cells.Item(1,1).Value = "Datos actualizados a 7 de Julio de 2020 a las 08:31"

# --- Updated per-country statistics ---
# Row 6 (India)
$ws.Cells.Item(6,2).Value = 720707
$ws.Cells.Item(6,3).Value = 361
$ws.Cells.Item(6,4).Value = 440192
$ws.Cells.Item(6,5).Value = 260337
$ws.Cells.Item(6,7).Value = 4
$ws.Cells.Item(6,8).Value = 20178

# Row 19 (Alemania)
$ws.Cells.Item(19,2).Value = 198064
$ws.Cells.Item(19,3).Value = 7
$ws.Cells.Item(19,5).Value = 6272

# Row 47 (Afganistan)
$ws.Cells.Item(47,2).Value = 33384
$ws.Cells.Item(47,3).Value = 194
$ws.Cells.Item(47,4).Value = 20179
$ws.Cells.Item(47,5).Value = 12285
$ws.Cells.Item(47,7).Value = 22
$ws.Cells.Item(47,8).Value = 920

# Row 49 (Israel)
$ws.Cells.Item(49,2).Value = 31186
$ws.Cells.Item(49,3).Value = 437
$ws.Cells.Item(49,4).Value = 18131
$ws.Cells.Item(49,5).Value = 12717
$ws.Cells.Item(49,7).Value = 4
$ws.Cells.Item(49,8).Value = 338

# Row 65 (now Camerun - fresh stats)
$ws.Cells.Item(65,2).Value = 14916
$ws.Cells.Item(65,3).Value = 2324
$ws.Cells.Item(65,4).Value = 11525
$ws.Cells.Item(65,5).Value = 3032
$ws.Cells.Item(65,7).Value = 46
$ws.Cells.Item(65,8).Value = 359

# Row 66 (now Marruecos - carries old Camerun-row's former numbers)
$ws.Cells.Item(66,2).Value = 14379
$ws.Cells.Item(66,3).Value = 0
$ws.Cells.Item(66,4).Value = 10173
$ws.Cells.Item(66,5).Value = 3969
$ws.Cells.Item(66,7).Value = 0
$ws.Cells.Item(66,8).Value = 237

# Row 67 (now Corea del Sur)
$ws.Cells.Item(67,2).Value = 13181
$ws.Cells.Item(67,3).Value = 44
$ws.Cells.Item(67,4).Value = 11914
$ws.Cells.Item(67,5).Value = 982
$ws.Cells.Item(67,7).Value = 1
$ws.Cells.Item(67,8).Value = 285

# Row 68 (now Dinamarca)
$ws.Cells.Item(68,2).Value = 12878
$ws.Cells.Item(68,4).Value = 11935
$ws.Cells.Item(68,5).Value = 336
$ws.Cells.Item(68,8).Value = 607

# Row 71 (Uzbekistan)
$ws.Cells.Item(71,2).Value = 10459
$ws.Cells.Item(71,3).Value = 97
$ws.Cells.Item(71,4).Value = 6690
$ws.Cells.Item(71,5).Value = 3732

# Row 78 (Australia)
$ws.Cells.Item(78,4).Value = 4785
$ws.Cells.Item(78,5).Value = 3019

# Row 112 (Sri Lanka)
$ws.Cells.Item(112,2).Value = 2078
$ws.Cells.Item(112,3).Value = 1
$ws.Cells.Item(112,5).Value = 150

# Row 141 (Georgia)
$ws.Cells.Item(141,2).Value = 958
$ws.Cells.Item(141,3).Value = 5
$ws.Cells.Item(141,4).Value = 838
$ws.Cells.Item(141,5).Value = 105
